# Auto-generated edit script: updates market-price-derived profit
# columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets to
# reflect a refreshed data pull. WVR is untouched.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1275.0978
$ws.Range("I15").Value = 1275.0978
$ws.Range("K15").Value = 3825.2934
$ws.Range("M15").Value = -3656.2934
$ws.Range("H17").Value = 2013.85
$ws.Range("J17").Value = 2013.85
$ws.Range("L17").Value = 6041.549999999999
$ws.Range("N17").Value = -6377.549999999999
$ws.Range("H96").Value = 745.8
$ws.Range("I96").Value = 723.3
$ws.Range("K96").Value = 2169.9
$ws.Range("M96").Value = -796.8999999999996
$ws.Range("H97").Value = 3768
$ws.Range("J97").Value = 3768
$ws.Range("L97").Value = 11304
$ws.Range("N97").Value = -12296
$ws.Range("H125").Value = 7939414
$ws.Range("I125").Value = 808.2
$ws.Range("K125").Value = 7273.8
$ws.Range("M125").Value = -4813.8
$ws.Range("H137").Value = 3255.465
$ws.Range("J137").Value = 5505.1
$ws.Range("L137").Value = 16515.3
$ws.Range("N137").Value = -21615.3
$ws.Range("H141").Value = 5846.074
$ws.Range("I141").Value = 5797.76
$ws.Range("K141").Value = 17393.28
$ws.Range("M141").Value = -12213.28

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4283.1514
$ws.Range("I32").Value = 3119.913
$ws.Range("J32").Value = 19571.428
$ws.Range("K32").Value = 3119.913
$ws.Range("L32").Value = 19571.428
$ws.Range("M32").Value = -2832.913
$ws.Range("N32").Value = -20145.428
$ws.Range("H64").Value = 100000
$ws.Range("J64").Value = 100000
$ws.Range("L64").Value = 100000
$ws.Range("N64").Value = -100496
$ws.Range("H67").Value = 100000
$ws.Range("J67").Value = 100000
$ws.Range("L67").Value = 100000
$ws.Range("N67").Value = -101716
$ws.Range("H110").Value = 187116.27
$ws.Range("I110").Value = 229215.5
$ws.Range("J110").Value = 1879.6
$ws.Range("K110").Value = 229215.5
$ws.Range("L110").Value = 1879.6
$ws.Range("M110").Value = -227170.5
$ws.Range("N110").Value = -5969.6
$ws.Range("H132").Value = 4328.185
$ws.Range("I132").Value = 3404.7046
$ws.Range("J132").Value = 8391.5
$ws.Range("K132").Value = 10214.1138
$ws.Range("L132").Value = 25174.5
$ws.Range("M132").Value = -7684.113799999999
$ws.Range("N132").Value = -30234.5

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H107").Value = 3338337.8
$ws.Range("J107").Value = 5005006.5
$ws.Range("L107").Value = 5005006.5
$ws.Range("N107").Value = -5008846.5

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 779.53845
$ws.Range("I22").Value = 463.2
$ws.Range("K22").Value = 463.2
$ws.Range("M22").Value = -113.2
$ws.Range("H31").Value = 46148.695
$ws.Range("I31").Value = 1892.6
$ws.Range("J31").Value = 80191.84
$ws.Range("K31").Value = 1892.6
$ws.Range("L31").Value = 80191.84
$ws.Range("M31").Value = -1597.6
$ws.Range("N31").Value = -80781.84
$ws.Range("H34").Value = 46148.695
$ws.Range("I34").Value = 1892.6
$ws.Range("J34").Value = 80191.84
$ws.Range("K34").Value = 1892.6
$ws.Range("L34").Value = 80191.84
$ws.Range("M34").Value = -1690.6
$ws.Range("N34").Value = -80595.84
$ws.Range("H50").Value = 30104.291
$ws.Range("J50").Value = 46799.582
$ws.Range("L50").Value = 46799.582
$ws.Range("N50").Value = -48049.582
$ws.Range("H122").Value = 1831.9231
$ws.Range("J122").Value = 3353.8
$ws.Range("L122").Value = 10061.4
$ws.Range("N122").Value = -14961.4
$ws.Range("H141").Value = 185797.23
$ws.Range("J141").Value = 183998.16
$ws.Range("L141").Value = 183998.16
$ws.Range("N141").Value = -194358.16

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 114970.72
$ws.Range("I107").Value = 1144.125
$ws.Range("J107").Value = 206032
$ws.Range("K107").Value = 3432.375
$ws.Range("L107").Value = 618096
$ws.Range("M107").Value = -1512.375
$ws.Range("N107").Value = -621936
$ws.Range("H129").Value = 64070.375
$ws.Range("I129").Value = 672.5714
$ws.Range("J129").Value = 113379.78
$ws.Range("K129").Value = 2017.7142
$ws.Range("L129").Value = 340139.34
$ws.Range("M129").Value = 2982.2858
$ws.Range("N129").Value = -350139.34
$ws.Range("H137").Value = 5855.857
$ws.Range("I137").Value = 3708.5557
$ws.Range("K137").Value = 11125.6671
$ws.Range("M137").Value = -6025.667099999999
$ws.Range("H138").Value = 3967.2
$ws.Range("I138").Value = 1991.5625
$ws.Range("J138").Value = 11869.75
$ws.Range("K138").Value = 5974.6875
$ws.Range("L138").Value = 35609.25
$ws.Range("M138").Value = -834.6875
$ws.Range("N138").Value = -45889.25

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 747363.5
$ws.Range("I80").Value = 594141.25
$ws.Range("J80").Value = 1007841.4
$ws.Range("K80").Value = 594141.25
$ws.Range("L80").Value = 1007841.4
$ws.Range("M80").Value = -593143.25
$ws.Range("N80").Value = -1009837.4
$ws.Range("H83").Value = 747363.5
$ws.Range("I83").Value = 594141.25
$ws.Range("J83").Value = 1007841.4
$ws.Range("K83").Value = 2970706.25
$ws.Range("L83").Value = 5039207
$ws.Range("M83").Value = -2965714.25
$ws.Range("N83").Value = -5049191
$ws.Range("H122").Value = 4972.75
$ws.Range("J122").Value = 5947.25
$ws.Range("L122").Value = 17841.75
$ws.Range("N122").Value = -22741.75

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7949.3335
$ws.Range("J7").Value = 7979.2
$ws.Range("L7").Value = 7979.2
$ws.Range("N7").Value = -8203.200000000001
$ws.Range("H18").Value = 35000
$ws.Range("J18").Value = 60000
$ws.Range("L18").Value = 60000
$ws.Range("N18").Value = -60344
$ws.Range("H20").Value = 618312.5
$ws.Range("I20").Value = 14300
$ws.Range("K20").Value = 14300
$ws.Range("M20").Value = -14074
$ws.Range("H30").Value = 5511.7144
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H40").Value = 337668
$ws.Range("I40").Value = 502502
$ws.Range("K40").Value = 502502
$ws.Range("M40").Value = -502366
$ws.Range("H46").Value = 2720.7856
$ws.Range("I46").Value = 2763.5454
$ws.Range("K46").Value = 2763.5454
$ws.Range("M46").Value = -2575.5454
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H68").Value = 126849.75
$ws.Range("I68").Value = 1979.6
$ws.Range("K68").Value = 1979.6
$ws.Range("M68").Value = -1230.6
$ws.Range("H70").Value = 10081.5
$ws.Range("J70").Value = 10081.5
$ws.Range("L70").Value = 10081.5
$ws.Range("N70").Value = -10621.5
$ws.Range("H71").Value = 126849.75
$ws.Range("I71").Value = 1979.6
$ws.Range("K71").Value = 9898
$ws.Range("M71").Value = -6154
$ws.Range("H73").Value = 10081.5
$ws.Range("J73").Value = 10081.5
$ws.Range("L73").Value = 10081.5
$ws.Range("N73").Value = -11953.5
$ws.Range("H122").Value = 1005439.6
$ws.Range("I122").Value = 1116599.5
$ws.Range("K122").Value = 3349798.5
$ws.Range("M122").Value = -3347348.5
$ws.Range("H126").Value = 7949.3335
$ws.Range("J126").Value = 7979.2
$ws.Range("L126").Value = 23937.6
$ws.Range("N126").Value = -28877.6
$ws.Range("H132").Value = 5987.2104
$ws.Range("I132").Value = 5007.5835
$ws.Range("K132").Value = 15022.7505
$ws.Range("M132").Value = -12492.7505

